$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as plain text in the source sheet (several
# look numeric, e.g. "345.47", and a few use "." as a thousands separator,
# e.g. "27.459.57"). Force each target cell to Text format first so Excel
# does not auto-convert the assigned string into a number.
$dCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.459.57"
$ws.Range("E2").Value = "  +6.53%  "
$ws.Range("D3").Value = "1.812.22"
$ws.Range("E3").Value = "  +6.36%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "345.47"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.3840"
$ws.Range("D8").Value = "50.25"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("D9").Value = "0.3519"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +12.78%  "
$ws.Range("D14").Value = "6.628"
$ws.Range("E14").Value = "  +7.14%  "
$ws.Range("D15").Value = "7.226"
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D16").Value = "1.811.23"
$ws.Range("E16").Value = "  +6.34%  "
$ws.Range("E17").Value = "  +5.71%  "
$ws.Range("D18").Value = "0.06761"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "86.87"
$ws.Range("E19").Value = "  +7.21%  "
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "17.85"
$ws.Range("E21").Value = "  +10.48%  "
$ws.Range("D22").Value = "6.531"
$ws.Range("E22").Value = "  +8.15%  "
$ws.Range("D23").Value = "13.19"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").Value = "27.468.08"
$ws.Range("E24").Value = "  +6.70%  "
$ws.Range("D25").Value = "2.470"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "2.682"
$ws.Range("E26").Value = "  +8.40%  "
$ws.Range("D27").Value = "22.13"
$ws.Range("E27").Value = "  +15.65%  "
$ws.Range("D28").Value = "1.504"
$ws.Range("E28").Value = "  +16.07%  "
$ws.Range("D29").Value = "154.29"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "2.013.36"
$ws.Range("E30").Value = "  +6.42%  "
$ws.Range("D31").Value = "137.01"
$ws.Range("E31").Value = "  +7.34%  "
$ws.Range("D32").Value = "6.391"
$ws.Range("E32").Value = "  +7.79%  "
$ws.Range("D33").Value = "4.083"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "13.93"
$ws.Range("E34").Value = "  +8.82%  "
$ws.Range("D35").Value = "0.08834"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").Value = "1.724"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "5.658"
$ws.Range("E37").Value = "  +6.56%  "
$ws.Range("D38").Value = "0.7101"
$ws.Range("E38").Value = "  +16.50%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2272"
$ws.Range("E39").Value = "  +7.44%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06546"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.02423"
$ws.Range("E41").Value = "  +7.90%  "
$ws.Range("D42").Value = "9.008"
$ws.Range("E42").Value = "  +5.58%  "
$ws.Range("D43").Value = "1.293"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "14.97"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "0.6622"
$ws.Range("E45").Value = "  +13.79%  "
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "3.985"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("D48").Value = "2.191"
$ws.Range("E48").Value = "  +9.65%  "
$ws.Range("D49").Value = "133.32"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").Value = "0.07374"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("D51").Value = "80.77"
$ws.Range("E51").Value = "  +5.79%  "

# Remove the temporary Text number format again so the cells end up with no
# explicit style (matching the original, unstyled data cells).
foreach ($addr in $dCells) {
    $ws.Range($addr).ClearFormats()
}
